$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Helper: replace the text of a whole paragraph (identified by its 1-based
# Paragraphs() index) with a sequence of text segments, and make sure each
# segment ends up in its own run (w:r) in the saved OOXML - mirroring the
# multi-run split that shows up in the target diff. Word normally coalesces
# adjacent same-formatted text into a single run, so to force genuine run
# boundaries we briefly flip Font.Bold on a sub-range and flip it back;
# that is enough to make the engine keep the runs distinct while leaving
# the visible formatting unchanged.
# -----------------------------------------------------------------------
function Set-ParagraphRuns($doc, [int]$paraIndex, [string[]]$segments) {
    $para = $doc.Paragraphs($paraIndex).Range
    $start = $para.Start
    $pEndExclMark = $para.End - 1

    $joined = [string]::Join("", $segments)
    $full = $doc.Range($start, $pEndExclMark)
    $full.Text = $joined

    # cumulative offsets where a new run should begin (skip the very first
    # segment, which starts at $start and needs no split point)
    $offsets = @()
    $cursor = $start
    for ($i = 0; $i -lt $segments.Length - 1; $i++) {
        $cursor = $cursor + $segments[$i].Length
        $offsets += $cursor
    }

    $newEnd = $start + $joined.Length

    # Apply splits back-to-front so each toggle only ever touches the
    # still-unsplit tail of the paragraph.
    for ($i = $offsets.Length - 1; $i -ge 0; $i--) {
        $pos = $offsets[$i]
        $r = $doc.Range($pos, $newEnd)
        $cur = $r.Font.Bold
        $r.Font.Bold = 1
        $r.Font.Bold = $cur
    }
}

# --- "Project Selection" paragraph ---------------------------------------
Set-ParagraphRuns $d 2 @(
    "We have selected project ",
    "1",
    " after the group discussion in the meeting. It was a web and ",
    "data analytics-based",
    " application."
)

# --- "Task Assign" paragraph ----------------------------------------------
Set-ParagraphRuns $d 5 @(
    "I was assigned ",
    "vueStoreFront",
    " application for installation on my local system."
)

# --- "Setup" paragraph -----------------------------------------------------
Set-ParagraphRuns $d 8 @(
    "The project is available in a code base of ",
    "GitHub",
    ". We divided this into 3 parts as was mentioned in the guidelines. Web, ",
    "and Data Analytics",
    ". I chose to work on ",
    "vueStoreFront",
    " application. It was a simple and very easy process to set up this."
)

# --- "I clone the ... application." bullet ---------------------------------
Set-ParagraphRuns $d 9 @(
    "I clone the ",
    "vueStoreFront",
    " application."
)

# --- "It contains the code of ..." bullet -----------------------------------
Set-ParagraphRuns $d 10 @(
    "It contains the code of ",
    "typescript."
)

# --- Remove the "As it isn't built directly..." bullet and the blank
#     paragraph right after it -------------------------------------------
$removeIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*As it isn*built directly*") {
        $removeIdx = $i
    }
}
if ($removeIdx -gt 0) {
    $p1 = $d.Paragraphs($removeIdx)
    $p2 = $d.Paragraphs($removeIdx + 1)
    $rng = $d.Range($p1.Range.Start, $p2.Range.End)
    $rng.Delete()
}
